$d = $word.ActiveDocument

# Insert point: the very end of the document body (after the last
# paragraph's text, before the final section break), so the existing
# content is preserved and the new paragraphs are appended after it.
$r = $d.Range($d.Content.End, $d.Content.End)

$rPr = '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'

$pPr1 = '<w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="240" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/>' + $rPr + '</w:pPr>'
$pPr2 = '<w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="2"/></w:numPr><w:spacing w:before="240" w:line="360" w:lineRule="auto"/><w:jc w:val="both"/>' + $rPr + '</w:pPr>'

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
$pPr1
<w:r>$rPr<w:t xml:space="preserve"> </w:t></w:r>
<w:r>$rPr<w:t>Trabalhando com fontes alternativas</w:t></w:r>
<w:r>$rPr<w:t>:</w:t></w:r>
</w:p>
<w:p>
$pPr2
<w:r>$rPr<w:t xml:space="preserve">Pegamos fontes open </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r>$rPr<w:t>source</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r>$rPr<w:t xml:space="preserve"> do google e colocamos no nosso código.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$r.InsertXML($xml)
